$d = $word.ActiveDocument

# 1. Remove the "Zdroje: " run that precedes the first hyperlink.
$zdroje = $d.Range(0, 8)
if ($zdroje.Text -eq "Zdroje: ") {
    $zdroje.Delete()
}

# 2. Remove the last two paragraphs ("Zvolil jsem flutter..." and
#    "Pro databazi pouziji PostgresSQL...") entirely, leaving the
#    bookmark-only paragraph behind as a single empty paragraph.
$count = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($count)
$pLast.Range.Delete()

$pPrev = $d.Paragraphs.Item($count - 1)
$pPrev.Range.Delete()

# 3. Strip the now-orphaned "_GoBack" bookmark so the remaining
#    paragraph becomes a clean empty paragraph.
$d.Bookmarks.Item("_GoBack").Delete()
